$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (within rich-text shared strings) ---
$ws.Range("A8").Value = "Volume 30   Number  9"
$ws.Range("C9").Value = "Report Covering the Week  2/27/2023  Through  3/5/2023"

# --- Simple numeric value updates (rows 14-30 crime stats) ---
$ws.Range("N14").Value = -66.666666666666
$ws.Range("N15").Value = -50
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = 9.090909090909
$ws.Range("F16").Value = 37
$ws.Range("G16").Value = 47
$ws.Range("H16").Value = -21.276595744680
$ws.Range("I16").Value = 95
$ws.Range("J16").Value = 84
$ws.Range("K16").Value = 13.095238095238
$ws.Range("L16").Value = 66.666666666666
$ws.Range("M16").Value = 196.875
$ws.Range("N16").Value = -77.647058823529
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 49
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 188.235294117647
$ws.Range("I17").Value = 105
$ws.Range("J17").Value = 43
$ws.Range("K17").Value = 144.186046511628
$ws.Range("L17").Value = 66.666666666666
$ws.Range("M17").Value = 169.230769230769
$ws.Range("N17").Value = -7.894736842105
$ws.Range("C18").Value = 12
$ws.Range("D18").Value = 14
$ws.Range("E18").Value = -14.285714285714
$ws.Range("F18").Value = 38
$ws.Range("G18").Value = 53
$ws.Range("H18").Value = -28.301886792452
$ws.Range("I18").Value = 91
$ws.Range("J18").Value = 101
$ws.Range("K18").Value = -9.900990099009
$ws.Range("L18").Value = 31.884057971014
$ws.Range("M18").Value = 30
$ws.Range("N18").Value = -82.398452611218
$ws.Range("C19").Value = 36
$ws.Range("D19").Value = 40
$ws.Range("E19").Value = -10
$ws.Range("F19").Value = 170
$ws.Range("G19").Value = 155
$ws.Range("H19").Value = 9.677419354838
$ws.Range("I19").Value = 407
$ws.Range("J19").Value = 325
$ws.Range("K19").Value = 25.230769230769
$ws.Range("L19").Value = 151.234567901235
$ws.Range("M19").Value = 2.518891687657
$ws.Range("N19").Value = -76.795895096921
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 9
$ws.Range("K20").Value = 28.571428571428
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 50
$ws.Range("N20").Value = -87.5
$ws.Range("C21").Value = 64
$ws.Range("E21").Value = -7.246376811594
$ws.Range("F21").Value = 300
$ws.Range("G21").Value = 275
$ws.Range("H21").Value = 9.090909090909
$ws.Range("I21").Value = 710
$ws.Range("J21").Value = 562
$ws.Range("K21").Value = 26.334519572953
$ws.Range("L21").Value = 95.592286501377
$ws.Range("M21").Value = 30.275229357798
$ws.Range("N21").Value = -75.424022152994
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 7
$ws.Range("E22").Value = -85.714285714285
$ws.Range("F22").Value = 11
$ws.Range("G22").Value = 16
$ws.Range("H22").Value = -31.25
$ws.Range("I22").Value = 30
$ws.Range("J22").Value = 31
$ws.Range("K22").Value = -3.225806451612
$ws.Range("L22").Value = -9.090909090909
$ws.Range("M22").Value = 7.142857142857
$ws.Range("C24").Value = 63
$ws.Range("D24").Value = 50
$ws.Range("E24").Value = 26
$ws.Range("F24").Value = 258
$ws.Range("G24").Value = 222
$ws.Range("H24").Value = 16.216216216216
$ws.Range("I24").Value = 622
$ws.Range("J24").Value = 439
$ws.Range("K24").Value = 41.685649202733
$ws.Range("L24").Value = 76.704545454545
$ws.Range("M24").Value = -19.010416666666
$ws.Range("C25").Value = 16
$ws.Range("E25").Value = -5.882352941176
$ws.Range("F25").Value = 64
$ws.Range("G25").Value = 66
$ws.Range("H25").Value = -3.030303030303
$ws.Range("I25").Value = 149
$ws.Range("J25").Value = 137
$ws.Range("K25").Value = 8.759124087591
$ws.Range("L25").Value = 19.2
$ws.Range("M25").Value = 49
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 5
$ws.Range("K26").Value = 150
$ws.Range("L26").Value = -16.666666666666
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 300
$ws.Range("F27").Value = 16
$ws.Range("G27").Value = 13
$ws.Range("H27").Value = 23.076923076923
$ws.Range("I27").Value = 37
$ws.Range("J27").Value = 29
$ws.Range("K27").Value = 27.586206896551
$ws.Range("L27").Value = 131.25
$ws.Range("G30").Value = 6
$ws.Range("H30").Value = -83.333333333333
$ws.Range("I30").Value = 3
$ws.Range("J30").Value = 6
$ws.Range("K30").Value = -50

# --- Cells that flip between numeric and text N/A-placeholder: copy content+style from a matching reference cell, not just a raw .Value= (which would attach the wrong style/number-format), then pin the exact value ---
$ws.Range("C23").Copy()
$ws.Range("D15").PasteSpecial(-4104)
$ws.Range("C23").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E23").Copy()
$ws.Range("E15").PasteSpecial(-4104)
$ws.Range("C23").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("C23").Copy()
$ws.Range("C20").PasteSpecial(-4104)
$ws.Range("C23").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C23").Copy()
$ws.Range("D26").PasteSpecial(-4104)
$ws.Range("C23").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E23").Copy()
$ws.Range("E26").PasteSpecial(-4104)
$ws.Range("C23").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("F30").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = 4
$ws.Range("H30").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100

$excel.CutCopyMode = $false
$ws.Range("A1").Select()
